{"js": "// The BEM lesson document gets a new paragraph explaining why BEM classes use\n// \"_\" / \"-\" separators. It is inserted right after the blank paragraph that\n// follows \"Elementy v bloku se NEZANO\u0158UJ\u00cd!! ...\" and right before \"Pro\u010d? \".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst newText =\n  \"D\u016fvod pro\u010d jsou _ _ nebo - - je proto\u017ee pokud je t\u0159\u00edda v\u00edceslovn\u00fd n\u00e1zev, odd\u011bluj\u00ed se \u2013 nebo _ , sp\u00ed\u0161 ne\u017e had\u00ed notace\";\n\n// Find the blank paragraph that sits right before the \"Pro\u010d? \" paragraph \u2013\n// that's the anchor point the new paragraph needs to follow.\nlet anchor = null;\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.trim() === \"\" && items[i + 1] && items[i + 1].text.trim() === \"Pro\u010d?\") {\n    anchor = items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not locate the blank paragraph preceding 'Pro\u010d?'.\");\n}\n\nanchor.insertParagraph(newText, \"After\");\nawait context.sync();\n", "ps1": "# The BEM lesson document gets a new paragraph explaining why BEM classes use\n# \"_\" / \"-\" separators. It is inserted right after the blank paragraph that\n# follows \"Elementy v bloku se NEZANO\u0158UJ\u00cd!! ...\" and right before \"Pro\u010d? \".\n$d = $word.ActiveDocument\n\n$newText = \"D\u016fvod pro\u010d jsou _ _ nebo - - je proto\u017ee pokud je t\u0159\u00edda v\u00edceslovn\u00fd n\u00e1zev, odd\u011bluj\u00ed se \u2013 nebo _ , sp\u00ed\u0161 ne\u017e had\u00ed notace\"\n\n# Locate the blank paragraph that immediately precedes the \"Pro\u010d?\" paragraph -\n# that's the anchor the new paragraph has to follow.\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq \"\" -and $i -lt $count) {\n        $next = $d.Paragraphs.Item($i + 1)\n        if ($next.Range.Text.Trim() -eq \"Pro\u010d?\") {\n            $targetIndex = $i\n            break\n        }\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the blank paragraph preceding 'Pro\u010d?'.\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n$target.Range.InsertParagraphAfter()\n$d.Paragraphs.Item($targetIndex + 1).Range.Text = $newText\n"}
